# Remove "name" from TextFileSequence (sequence's file-source sheet) and
# reorder Primer's header row so "name" moves to the front.

$wb = $excel.ActiveWorkbook

# --- TextFileSequence: drop the trailing "name" column (G1) -----------------
$wsSeq = $wb.Worksheets.Item("TextFileSequence")
$wsSeq.Range("G1").ClearContents()

# --- Primer: reorder header row from (sequence, id, type, name) -------------
#     to (name, sequence, id, type)
$wsPrimer = $wb.Worksheets.Item("Primer")
$wsPrimer.Range("A1").Value = "name"
$wsPrimer.Range("B1").Value = "sequence"
$wsPrimer.Range("C1").Value = "id"
$wsPrimer.Range("D1").Value = "type"
